$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# The handback transform for the ffefa583-... file (row 3 on every sheet)
# failed, so its shared "Status" text changes everywhere it is shown.
$wsOverview.Range("E3").Value = "Handback transform failed"
$wsOverview.Range("F3").Value = "Handback transform failed"
$wsZhCn.Range("C3").Value = "Handback transform failed"
$wsDeDe.Range("C3").Value = "Handback transform failed"

# Record the handback/handoff file name mismatch error in the "Error Detail"
# column (P) of the ffefa583 row on each localized-language sheet.
$wsZhCn.Range("P3").Value = "Handback file name: ocwjaa00.xgr is different with handoff file name: ffefa583-e5dc-4375-82d4-3555298fb638.789dc666e32550f0306b93ca3f2a595a4e10d327.zh-cn."
$wsDeDe.Range("P3").Value = "Handback file name: ocwjaa00.xgr is different with handoff file name: ffefa583-e5dc-4375-82d4-3555298fb638.789dc666e32550f0306b93ca3f2a595a4e10d327.de-de."

# Widen column P so the new error message is readable. ColumnWidth is stored
# internally with a fixed +5/6 character padding baked in by this host, so
# back that off here to land on a stored column width of exactly 40.
$wsZhCn.Range("P1").ColumnWidth = 39.166666666666664
$wsDeDe.Range("P1").ColumnWidth = 39.166666666666664
